# Auto-generated Excel COM-interop edit script
# Applies the changes described in the commit diff to global_sdg_indicators.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Last update" timestamp note (shared string reused via cell A472) ---
$ws.Range("A472").Value = "Last update: 25-10-2022, 11:25"

# --- Unfreeze panes and set plain selection on A1:B1 (matches new sheetView state) ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A1:B1").Select()

# --- Copy number formatting (style) from same-row reference cell into previously-blank cells ---
# so that newly populated cells adopt the same numeric style as their row neighbours,
# exactly as Excel does when a user fills in an empty cell next to existing data.
$ws.Range("F376").Copy()
$ws.Range("P376").PasteSpecial(-4122)
$ws.Range("F377").Copy()
$ws.Range("P377").PasteSpecial(-4122)
$ws.Range("F378").Copy()
$ws.Range("P378").PasteSpecial(-4122)

$ws.Range("O33").Copy()
$ws.Range("P33").PasteSpecial(-4122)
$ws.Range("O34").Copy()
$ws.Range("P34").PasteSpecial(-4122)
$ws.Range("O35").Copy()
$ws.Range("P35").PasteSpecial(-4122)
$ws.Range("O36").Copy()
$ws.Range("P36").PasteSpecial(-4122)
$ws.Range("O37").Copy()
$ws.Range("P37").PasteSpecial(-4122)
$ws.Range("O38").Copy()
$ws.Range("P38").PasteSpecial(-4122)
$ws.Range("O39").Copy()
$ws.Range("P39").PasteSpecial(-4122)
$ws.Range("O204").Copy()
$ws.Range("P204").PasteSpecial(-4122)

$ws.Range("P75").Copy()
$ws.Range("Q75").PasteSpecial(-4122)
$ws.Range("P113").Copy()
$ws.Range("Q113").PasteSpecial(-4122)
$ws.Range("P114").Copy()
$ws.Range("Q114").PasteSpecial(-4122)
$ws.Range("P115").Copy()
$ws.Range("Q115").PasteSpecial(-4122)
$ws.Range("P116").Copy()
$ws.Range("Q116").PasteSpecial(-4122)
$ws.Range("P117").Copy()
$ws.Range("Q117").PasteSpecial(-4122)
$ws.Range("P119").Copy()
$ws.Range("Q119").PasteSpecial(-4122)
$ws.Range("P120").Copy()
$ws.Range("Q120").PasteSpecial(-4122)
$ws.Range("P121").Copy()
$ws.Range("Q121").PasteSpecial(-4122)
$ws.Range("P122").Copy()
$ws.Range("Q122").PasteSpecial(-4122)
$ws.Range("P123").Copy()
$ws.Range("Q123").PasteSpecial(-4122)
$ws.Range("P199").Copy()
$ws.Range("Q199").PasteSpecial(-4122)
$ws.Range("P200").Copy()
$ws.Range("Q200").PasteSpecial(-4122)
$ws.Range("P214").Copy()
$ws.Range("Q214").PasteSpecial(-4122)
$ws.Range("P257").Copy()
$ws.Range("Q257").PasteSpecial(-4122)
$ws.Range("P368").Copy()
$ws.Range("Q368").PasteSpecial(-4122)
$ws.Range("P369").Copy()
$ws.Range("Q369").PasteSpecial(-4122)
$ws.Range("P370").Copy()
$ws.Range("Q370").PasteSpecial(-4122)
$ws.Range("P371").Copy()
$ws.Range("Q371").PasteSpecial(-4122)
$ws.Range("P454").Copy()
$ws.Range("Q454").PasteSpecial(-4122)

$ws.Range("P22").Copy()
$ws.Range("Q22").PasteSpecial(-4122)
$ws.Range("P23").Copy()
$ws.Range("Q23").PasteSpecial(-4122)
$ws.Range("P256").Copy()
$ws.Range("Q256").PasteSpecial(-4122)
$ws.Range("P334").Copy()
$ws.Range("Q334").PasteSpecial(-4122)
$ws.Range("P335").Copy()
$ws.Range("Q335").PasteSpecial(-4122)
$ws.Range("P336").Copy()
$ws.Range("Q336").PasteSpecial(-4122)
$ws.Range("P337").Copy()
$ws.Range("Q337").PasteSpecial(-4122)
$ws.Range("P338").Copy()
$ws.Range("Q338").PasteSpecial(-4122)
$ws.Range("P339").Copy()
$ws.Range("Q339").PasteSpecial(-4122)
$ws.Range("P364").Copy()
$ws.Range("Q364").PasteSpecial(-4122)
$ws.Range("P365").Copy()
$ws.Range("Q365").PasteSpecial(-4122)
$ws.Range("P366").Copy()
$ws.Range("Q366").PasteSpecial(-4122)
$ws.Range("P367").Copy()
$ws.Range("Q367").PasteSpecial(-4122)
$ws.Range("P383").Copy()
$ws.Range("Q383").PasteSpecial(-4122)
$ws.Range("P384").Copy()
$ws.Range("Q384").PasteSpecial(-4122)
$ws.Range("P446").Copy()
$ws.Range("Q446").PasteSpecial(-4122)

$ws.Range("P156").Copy()
$ws.Range("Q156").PasteSpecial(-4122)
$ws.Range("P267").Copy()
$ws.Range("Q267").PasteSpecial(-4122)
$ws.Range("P268").Copy()
$ws.Range("Q268").PasteSpecial(-4122)
$ws.Range("P269").Copy()
$ws.Range("Q269").PasteSpecial(-4122)
$ws.Range("P270").Copy()
$ws.Range("Q270").PasteSpecial(-4122)
$ws.Range("P271").Copy()
$ws.Range("Q271").PasteSpecial(-4122)
$ws.Range("P272").Copy()
$ws.Range("Q272").PasteSpecial(-4122)
$ws.Range("P273").Copy()
$ws.Range("Q273").PasteSpecial(-4122)
$ws.Range("P274").Copy()
$ws.Range("Q274").PasteSpecial(-4122)
$ws.Range("P275").Copy()
$ws.Range("Q275").PasteSpecial(-4122)
$ws.Range("P276").Copy()
$ws.Range("Q276").PasteSpecial(-4122)
$ws.Range("P277").Copy()
$ws.Range("Q277").PasteSpecial(-4122)
$ws.Range("P389").Copy()
$ws.Range("Q389").PasteSpecial(-4122)
$ws.Range("P390").Copy()
$ws.Range("Q390").PasteSpecial(-4122)
$ws.Range("P391").Copy()
$ws.Range("Q391").PasteSpecial(-4122)
$ws.Range("P392").Copy()
$ws.Range("Q392").PasteSpecial(-4122)
$ws.Range("P405").Copy()
$ws.Range("Q405").PasteSpecial(-4122)
$ws.Range("P409").Copy()
$ws.Range("Q409").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Set new cell values ---
$ws.Range("Q22").Value = 2.25
$ws.Range("Q23").Value = 0.68
$ws.Range("P33").Value = 2.5
$ws.Range("P34").Value = 7.4
$ws.Range("P35").Value = 8.4
$ws.Range("P36").Value = 5.3
$ws.Range("P37").Value = 0.9
$ws.Range("O38").Value = 0.5
$ws.Range("P38").Value = 1.1000000000000001
$ws.Range("P39").Value = 0.9
$ws.Range("P75").Value = 10.9
$ws.Range("Q75").Value = 11
$ws.Range("Q113").Value = 28.4
$ws.Range("Q114").Value = 51.4
$ws.Range("Q115").Value = 35.799999999999997
$ws.Range("Q116").Value = 33.9
$ws.Range("Q117").Value = 68.3
$ws.Range("Q119").Value = 55.2
$ws.Range("Q120").Value = 87.5
$ws.Range("Q121").Value = 69.8
$ws.Range("Q122").Value = 72.7
$ws.Range("Q123").Value = 91.5
$ws.Range("Q156").Value = 117
$ws.Range("Q199").Value = 62.4
$ws.Range("Q200").Value = 71.900000000000006
$ws.Range("P204").Value = 53.4
$ws.Range("F214").Value = 3.3
$ws.Range("G214").Value = 5
$ws.Range("H214").Value = 1.5
$ws.Range("I214").Value = 0.9
$ws.Range("J214").Value = 3.9
$ws.Range("K214").Value = 4.5
$ws.Range("L214").Value = 3
$ws.Range("M214").Value = 5.2
$ws.Range("N214").Value = 6
$ws.Range("O214").Value = 4.5
$ws.Range("P214").Value = -1.5
$ws.Range("Q214").Value = 7.4
$ws.Range("F215").Value = 13.8
$ws.Range("H215").Value = 4.0999999999999996
$ws.Range("I215").Value = 3
$ws.Range("J215").Value = 1.7
$ws.Range("K215").Value = 3.1
$ws.Range("L215").Value = 3.1
$ws.Range("M215").Value = 3.4
$ws.Range("N215").Value = 7.2
$ws.Range("O215").Value = 5.7
$ws.Range("P215").Value = 1.3
$ws.Range("Q215").Value = 8.5
$ws.Range("F256").Value = 0.17
$ws.Range("G256").Value = 0.2
$ws.Range("H256").Value = 0.16
$ws.Range("I256").Value = 0.14000000000000001
$ws.Range("J256").Value = 0.11
$ws.Range("K256").Value = 0.14000000000000001
$ws.Range("L256").Value = 0.11
$ws.Range("M256").Value = 0.16
$ws.Range("N256").Value = 0.15
$ws.Range("O256").Value = 0.17
$ws.Range("P256").Value = 0.07000000000000001
$ws.Range("Q256").Value = 0.1
$ws.Range("Q257").Value = 78.8
$ws.Range("Q267").Value = 244874
$ws.Range("Q268").Value = 168619
$ws.Range("Q269").Value = 674
$ws.Range("Q270").Value = 986
$ws.Range("Q271").Value = 6996
$ws.Range("Q272").Value = 237915
$ws.Range("P273").Value = 1919193
$ws.Range("Q273").Value = 1952465
$ws.Range("Q274").Value = 49855
$ws.Range("Q275").Value = 9587
$ws.Range("Q276").Value = 3465
$ws.Range("Q277").Value = 91
$ws.Range("F278").Value = 16.3
$ws.Range("G278").Value = 16.3
$ws.Range("H278").Value = 16.7
$ws.Range("I278").Value = 15.7
$ws.Range("J278").Value = 17.100000000000001
$ws.Range("K278").Value = 17.899999999999999
$ws.Range("L278").Value = 18.5
$ws.Range("M278").Value = 17.100000000000001
$ws.Range("N278").Value = 16.8
$ws.Range("O278").Value = 16.899999999999999
$ws.Range("P278").Value = 16.5
$ws.Range("Q278").Value = 16.7
$ws.Range("Q334").Value = 2.25
$ws.Range("Q335").Value = 0.68
$ws.Range("Q336").Value = 26.92
$ws.Range("Q337").Value = 13.34
$ws.Range("Q338").Value = 21.01
$ws.Range("Q339").Value = 38.729999999999997
$ws.Range("P364").Value = 98.41
$ws.Range("Q364").Value = 98.52
$ws.Range("Q365").Value = 100
$ws.Range("P366").Value = 95
$ws.Range("Q366").Value = 95.35
$ws.Range("P367").Value = 99.29
$ws.Range("Q367").Value = 99.32
$ws.Range("Q368").Value = 99.9
$ws.Range("Q369").Value = 100
$ws.Range("Q370").Value = 99.8
$ws.Range("Q371").Value = 99.9
$ws.Range("P376").Value = 59
$ws.Range("F377").Value = 36
$ws.Range("H377").Value = 37
$ws.Range("J377").Value = 35
$ws.Range("L377").Value = 34
$ws.Range("N377").Value = 37
$ws.Range("P377").Value = 36
$ws.Range("F378").Value = 42
$ws.Range("H378").Value = 44
$ws.Range("J378").Value = 46
$ws.Range("L378").Value = 44
$ws.Range("N378").Value = 44
$ws.Range("P378").Value = 42
$ws.Range("Q383").Value = 2.25
$ws.Range("Q384").Value = 0.68
$ws.Range("Q389").Value = 67
$ws.Range("Q390").Value = 17
$ws.Range("Q391").Value = 17
$ws.Range("Q392").Value = 33
$ws.Range("Q405").Value = 20928
$ws.Range("Q409").Value = 20928
$ws.Range("Q446").Value = 0.15
$ws.Range("Q454").Value = 85.4
$ws.Range("F457").Value = 37240
$ws.Range("G457").Value = 40327
$ws.Range("H457").Value = 41852
$ws.Range("I457").Value = 42339
$ws.Range("J457").Value = 44189
$ws.Range("K457").Value = 46768
$ws.Range("L457").Value = 48227
$ws.Range("M457").Value = 51606
$ws.Range("N457").Value = 55359
$ws.Range("O457").Value = 59618
$ws.Range("P457").Value = 61231
$ws.Range("Q457").Value = 69069
$ws.Range("F458").Value = 103.4
$ws.Range("G458").Value = 105
$ws.Range("H458").Value = 101.5
$ws.Range("I458").Value = 100.9
$ws.Range("J458").Value = 103.8
$ws.Range("K458").Value = 104.4
$ws.Range("L458").Value = 103
$ws.Range("M458").Value = 105.1
$ws.Range("N458").Value = 105.9
$ws.Range("O458").Value = 104.5
$ws.Range("P458").Value = 98
$ws.Range("Q458").Value = 106.8
$ws.Range("F460").Value = 19.7
$ws.Range("G460").Value = 20.5
$ws.Range("H460").Value = 19.600000000000001
$ws.Range("I460").Value = 19
$ws.Range("J460").Value = 20.100000000000001
$ws.Range("K460").Value = 20.399999999999999
$ws.Range("L460").Value = 18.5
$ws.Range("M460").Value = 17.600000000000001
$ws.Range("N460").Value = 18.7
$ws.Range("O460").Value = 18.899999999999999
$ws.Range("P460").Value = 18.3
$ws.Range("Q460").Value = 17
$ws.Range("N461").Value = 48.7
$ws.Range("O461").Value = 45.7
$ws.Range("P461").Value = 57.2
$ws.Range("Q462").Value = -1.8
